$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44659
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 350
$ws.Range("L2").Value = 380
$ws.Range("M2").Value = 365
$ws.Range("P2").Value = 365
$ws.Range("D3").Value = 44251
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 250
$ws.Range("L3").Value = 280
$ws.Range("M3").Value = 265
$ws.Range("P3").Value = 265
$ws.Range("D4").Value = 44650
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 325
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = 338
$ws.Range("P4").Value = 338
$ws.Range("D5").Value = 44566
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 320
$ws.Range("M5").Value = 310
$ws.Range("O5").Value = "Perú"
$ws.Range("P5").Value = 310
$ws.Range("D6").Value = 44243
$ws.Range("J6").Value = 1200
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 320
$ws.Range("M6").Value = 310
$ws.Range("P6").Value = 310
$ws.Range("D7").Value = 44243
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 800
$ws.Range("D8").Value = 44575
$ws.Range("J8").Value = 1200
$ws.Range("K8").Value = 380
$ws.Range("L8").Value = 400
$ws.Range("M8").Value = 390
$ws.Range("P8").Value = 390
$ws.Range("D9").Value = 44166
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1700
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 530
$ws.Range("M9").Value = 515
$ws.Range("P9").Value = 515
$ws.Range("D10").Value = 44523
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 550
$ws.Range("L10").Value = 580
$ws.Range("M10").Value = 565
$ws.Range("P10").Value = 565
$ws.Range("D11").Value = 44589
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 900
$ws.Range("K11").Value = 325
$ws.Range("L11").Value = 350
$ws.Range("M11").Value = 338
$ws.Range("P11").Value = 338
$ws.Range("D12").Value = 44214
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 450
$ws.Range("M12").Value = 425
$ws.Range("P12").Value = 425
$ws.Range("D13").Value = 44224
$ws.Range("J13").Value = 1200
$ws.Range("K13").Value = 230
$ws.Range("L13").Value = 250
$ws.Range("M13").Value = 240
$ws.Range("P13").Value = 240
$ws.Range("D14").Value = 44224
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 200
$ws.Range("L14").Value = 230
$ws.Range("M14").Value = 215
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 215
$ws.Range("D15").Value = 44160
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 550
$ws.Range("M15").Value = 525
$ws.Range("P15").Value = 525
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1300
$ws.Range("K16").Value = 350
$ws.Range("L16").Value = 380
$ws.Range("M16").Value = 365
$ws.Range("P16").Value = 365
$ws.Range("D17").Value = 44602
$ws.Range("J17").Value = 900
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 330
$ws.Range("M17").Value = 315
$ws.Range("O17").Value = "Perú"
$ws.Range("P17").Value = 315
$ws.Range("D18").Value = 44530
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 1300
$ws.Range("K18").Value = 450
$ws.Range("L18").Value = 480
$ws.Range("M18").Value = 465
$ws.Range("P18").Value = 465
$ws.Range("D19").Value = 44586
$ws.Range("I19").Value = "Tercera"
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 330
$ws.Range("L19").Value = 350
$ws.Range("M19").Value = 340
$ws.Range("O19").Value = "Región de Arica y Parinacota"
$ws.Range("P19").Value = 340
$ws.Range("D20").Value = 44603
$ws.Range("I20").Value = "Tercera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 280
$ws.Range("L20").Value = 300
$ws.Range("M20").Value = 290
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 290
$ws.Range("D21").Value = 44201
$ws.Range("J21").Value = 1800
$ws.Range("K21").Value = 250
$ws.Range("L21").Value = 270
$ws.Range("M21").Value = 260
$ws.Range("P21").Value = 260
$ws.Range("D22").Value = 44580
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 380
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 390
$ws.Range("O22").Value = "Región Metropolitana"
$ws.Range("P22").Value = 390
$ws.Range("D23").Value = 44168
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 1700
$ws.Range("K23").Value = 430
$ws.Range("L23").Value = 450
$ws.Range("M23").Value = 440
$ws.Range("O23").Value = "Perú"
$ws.Range("P23").Value = 440
$ws.Range("D24").Value = 44217
$ws.Range("J24").Value = 1600
$ws.Range("K24").Value = 300
$ws.Range("L24").Value = 350
$ws.Range("M24").Value = 325
$ws.Range("P24").Value = 325
$ws.Range("D25").Value = 44176
$ws.Range("J25").Value = 1300
$ws.Range("K25").Value = 350
$ws.Range("L25").Value = 400
$ws.Range("M25").Value = 375
$ws.Range("P25").Value = 375
$ws.Range("D26").Value = 44172
$ws.Range("J26").Value = 1600
$ws.Range("K26").Value = 400
$ws.Range("L26").Value = 420
$ws.Range("M26").Value = 410
$ws.Range("P26").Value = 410
$ws.Range("D27").Value = 44202
$ws.Range("J27").Value = 1300
$ws.Range("K27").Value = 230
$ws.Range("L27").Value = 250
$ws.Range("M27").Value = 240
$ws.Range("P27").Value = 240
$ws.Range("D28").Value = 44162
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 500
$ws.Range("L28").Value = 550
$ws.Range("M28").Value = 525
$ws.Range("P28").Value = 525
$ws.Range("D29").Value = 44162
$ws.Range("I29").Value = "Segunda"
$ws.Range("K29").Value = 500
$ws.Range("L29").Value = 550
$ws.Range("M29").Value = 525
$ws.Range("P29").Value = 525
$ws.Range("D30").Value = 44175
$ws.Range("I30").Value = "Segunda"
$ws.Range("K30").Value = 400
$ws.Range("L30").Value = 430
$ws.Range("M30").Value = 415
$ws.Range("O30").Value = "Perú"
$ws.Range("P30").Value = 415
$ws.Range("D31").Value = 44253
$ws.Range("I31").Value = "Segunda"
$ws.Range("J31").Value = 1200
$ws.Range("K31").Value = 270
$ws.Range("L31").Value = 280
$ws.Range("M31").Value = 275
$ws.Range("P31").Value = 275
$ws.Range("D32").Value = 44301
$ws.Range("J32").Value = 900
$ws.Range("K32").Value = 280
$ws.Range("L32").Value = 300
$ws.Range("M32").Value = 290
$ws.Range("P32").Value = 290
$ws.Range("D34").Value = 44231
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 180
$ws.Range("L34").Value = 200
$ws.Range("M34").Value = 190
$ws.Range("O34").Value = "Región de Arica y Parinacota"
$ws.Range("P34").Value = 190
$ws.Range("D35").Value = 44547
$ws.Range("J35").Value = 1200
$ws.Range("L35").Value = 370
$ws.Range("M35").Value = 360
$ws.Range("P35").Value = 360
$ws.Range("D36").Value = 44609
$ws.Range("K36").Value = 280
$ws.Range("L36").Value = 300
$ws.Range("M36").Value = 290
$ws.Range("P36").Value = 290
